# Update the quiz data on the active sheet with the new question set
# (lyrics/trivia about 雙截棍 "Shuang Jie Gun") replacing the old
# placeholder math / programming-language quiz questions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (single-choice question)
$ws.Range("B2").Value = "雙節棍中「快使用...」是什麼？"
$ws.Range("C2").Value = "雙截棍,刀,劍,槍,棒"
$ws.Range("D2").Value = "雙截棍"

# Row 3 (single-choice question)
$ws.Range("B3").Value = "歌詞提到的「夜市」是在什麼城市？"
$ws.Range("C3").Value = "台北,高雄,台中,台南,新竹"
$ws.Range("D3").Value = "台北"

# Row 4 (multiple-choice question)
$ws.Range("B4").Value = "哪些是歌詞中提到的武器？"
$ws.Range("C4").Value = "雙截棍,彈弓,匕首,長劍"
$ws.Range("D4").Value = "雙截棍,匕首"

# Row 5 (multiple-choice question)
$ws.Range("B5").Value = "哪些是歌詞提到的情景？"
$ws.Range("C5").Value = "夜市,擂台,公園,跑馬場"
$ws.Range("D5").Value = "夜市,擂台"

# Row 6 (fill-in-the-blank question)
$ws.Range("B6").Value = "「___ 使用雙截棍，哼哼哈嘿！」"
$ws.Range("D6").Value = "快"

# Match the selection shown on the refreshed sheet (whole used range selected)
$ws.Range("A1:D6").Select()
